# "Atu PPT e CPU 128bits 13092022"
#
# 1) Slide 2  - title placeholder: "Aula 01 e 02" -> "Aula 00"
# 2) Slide 25 - reference list: the hyperlink text was split across two
#    runs ("h" + "ttp://www.proedu...") - rejoin it into a single run
#    with the full, correctly spelled URL.

$p = $ppt.ActivePresentation

# --- 1) Slide 2: fix the course title -------------------------------------
$slide2 = $p.Slides.Item(2)
$titleShape = $slide2.Shapes.Item(4)          # "Título 2" (title placeholder)
$titleRange = $titleShape.TextFrame.TextRange

$oldTitle = "Aula 01 e 02"
$newTitle = "Aula 00"

$titleText = $titleRange.Text
$titleStart = $titleText.IndexOf($oldTitle) + 1   # 1-based start for Characters()
$titleRun = $titleRange.Characters($titleStart, $oldTitle.Length)
$titleRun.Text = $newTitle

# --- 2) Slide 25: repair the broken hyperlink run --------------------------
$slide25 = $p.Slides.Item(25)
$refShape = $slide25.Shapes.Item(2)            # "Text Placeholder 2"
$refRange = $refShape.TextFrame.TextRange

$fullUrl = "http://www.proedu.rnp.br/bitstream/handle/123456789/697/Arquitetura_de_Computadores_web.pdf?sequence=3&isAllowed=y"

$refText = $refRange.Text
$urlStart = $refText.IndexOf($fullUrl.Substring(1)) # 0-based idx of "ttp://..." == 1-based start of "http://..."
$urlRun = $refRange.Characters($urlStart, $fullUrl.Length)
$urlRun.Text = $fullUrl
